$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 672.38464
$ws.Range("I4").Value = 272.66666
$ws.Range("K4").Value = 272.66666
$ws.Range("M4").Value = -158.66666
$ws.Range("H28").Value = 53934.95
$ws.Range("I28").Value = 67421.60000000001
$ws.Range("J28").Value = 3360
$ws.Range("K28").Value = 67421.60000000001
$ws.Range("L28").Value = 3360
$ws.Range("M28").Value = -66936.60000000001
$ws.Range("N28").Value = -4330
$ws.Range("H70").Value = 1733.3334
$ws.Range("I70").Value = 1500
$ws.Range("J70").Value = 1780
$ws.Range("K70").Value = 4500
$ws.Range("L70").Value = 5340
$ws.Range("M70").Value = -4230
$ws.Range("N70").Value = -5880
$ws.Range("H73").Value = 1733.3334
$ws.Range("I73").Value = 1500
$ws.Range("J73").Value = 1780
$ws.Range("K73").Value = 4500
$ws.Range("L73").Value = 5340
$ws.Range("M73").Value = -3564
$ws.Range("N73").Value = -7212
$ws.Range("H80").Value = 14021.333
$ws.Range("I80").Value = 16741.715
$ws.Range("J80").Value = 4500
$ws.Range("K80").Value = 50225.145
$ws.Range("L80").Value = 13500
$ws.Range("M80").Value = -49227.145
$ws.Range("N80").Value = -15496
$ws.Range("H83").Value = 14021.333
$ws.Range("I83").Value = 16741.715
$ws.Range("J83").Value = 4500
$ws.Range("K83").Value = 150675.435
$ws.Range("L83").Value = 40500
$ws.Range("M83").Value = -145683.435
$ws.Range("N83").Value = -50484
$ws.Range("H86").Value = 20449.934
$ws.Range("I86").Value = 9723.25
$ws.Range("K86").Value = 9723.25
$ws.Range("M86").Value = -8600.25
$ws.Range("H89").Value = 20449.934
$ws.Range("I89").Value = 9723.25
$ws.Range("K89").Value = 48616.25
$ws.Range("M89").Value = -43000.25
$ws.Range("H96").Value = 3486.1
$ws.Range("J96").Value = 12000
$ws.Range("L96").Value = 36000
$ws.Range("N96").Value = -38746
$ws.Range("H100").Value = 6025.909
$ws.Range("I100").Value = 7785.625
$ws.Range("J100").Value = 1333.3334
$ws.Range("K100").Value = 7785.625
$ws.Range("L100").Value = 1333.3334
$ws.Range("M100").Value = -7244.625
$ws.Range("N100").Value = -2415.3334
$ws.Range("H112").Value = 31332.854
$ws.Range("I112").Value = 112843.664
$ws.Range("J112").Value = 1988.96
$ws.Range("K112").Value = 338530.992
$ws.Range("L112").Value = 5966.88
$ws.Range("M112").Value = -337422.992
$ws.Range("N112").Value = -8182.88
$ws.Range("H116").Value = 30761086
$ws.Range("I116").Value = 31377350
$ws.Range("K116").Value = 31377350
$ws.Range("M116").Value = -31373908
$ws.Range("H132").Value = 2528.6123
$ws.Range("J132").Value = 2679.5
$ws.Range("L132").Value = 8038.5
$ws.Range("N132").Value = -13098.5
$ws.Range("H141").Value = 3144.4119
$ws.Range("I141").Value = 2868.7144
$ws.Range("K141").Value = 8606.143199999999
$ws.Range("M141").Value = -3426.143199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1484.2333
$ws.Range("I32").Value = 1466.4482
$ws.Range("K32").Value = 1466.4482
$ws.Range("M32").Value = -1179.4482
$ws.Range("H45").Value = 2254.3333
$ws.Range("I45").Value = 2069.8572
$ws.Range("K45").Value = 2069.8572
$ws.Range("M45").Value = -1692.8572
$ws.Range("H97").Value = 856.3333
$ws.Range("I97").Value = 662.0833
$ws.Range("J97").Value = 1050.5834
$ws.Range("K97").Value = 662.0833
$ws.Range("L97").Value = 1050.5834
$ws.Range("M97").Value = -166.0833
$ws.Range("N97").Value = -2042.5834
$ws.Range("H102").Value = 4827288
$ws.Range("I102").Value = 6993612.5
$ws.Range("K102").Value = 6993612.5
$ws.Range("M102").Value = -6991990.5
$ws.Range("H110").Value = 62564310
$ws.Range("I110").Value = 83376590
$ws.Range("K110").Value = 83376590
$ws.Range("M110").Value = -83374545

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 23812184
$ws.Range("I86").Value = 33336172
$ws.Range("K86").Value = 33336172
$ws.Range("M86").Value = -33335049
$ws.Range("H89").Value = 23812184
$ws.Range("I89").Value = 33336172
$ws.Range("K89").Value = 166680860
$ws.Range("M89").Value = -166675244
$ws.Range("H94").Value = 3566.9092
$ws.Range("I94").Value = 3859.5625
$ws.Range("J94").Value = 2786.5
$ws.Range("K94").Value = 3859.5625
$ws.Range("L94").Value = 2786.5
$ws.Range("M94").Value = -3408.5625
$ws.Range("N94").Value = -3688.5
$ws.Range("H99").Value = 940.8125
$ws.Range("I99").Value = 911.8
$ws.Range("K99").Value = 911.8
$ws.Range("M99").Value = 586.2
$ws.Range("H134").Value = 3232.5454
$ws.Range("I134").Value = 3125.111
$ws.Range("K134").Value = 9375.332999999999
$ws.Range("M134").Value = -6840.332999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1026
$ws.Range("I16").Value = 1080.3334
$ws.Range("J16").Value = 700
$ws.Range("K16").Value = 1080.3334
$ws.Range("L16").Value = 700
$ws.Range("M16").Value = -793.3334
$ws.Range("N16").Value = -1274
$ws.Range("H62").Value = 142862990
$ws.Range("I62").Value = 7750
$ws.Range("J62").Value = 333336670
$ws.Range("K62").Value = 7750
$ws.Range("L62").Value = 333336670
$ws.Range("M62").Value = -7126
$ws.Range("N62").Value = -333337918
$ws.Range("H65").Value = 142862990
$ws.Range("I65").Value = 7750
$ws.Range("J65").Value = 333336670
$ws.Range("K65").Value = 38750
$ws.Range("L65").Value = 1666683350
$ws.Range("M65").Value = -35630
$ws.Range("N65").Value = -1666689590
$ws.Range("H113").Value = 1026
$ws.Range("I113").Value = 1080.3334
$ws.Range("J113").Value = 700
$ws.Range("K113").Value = 1080.3334
$ws.Range("L113").Value = 700
$ws.Range("M113").Value = 1089.6666
$ws.Range("N113").Value = -5040

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 499.14285
$ws.Range("I17").Value = 452.83334
$ws.Range("J17").Value = 777
$ws.Range("K17").Value = 1358.50002
$ws.Range("L17").Value = 2331
$ws.Range("M17").Value = -1189.50002
$ws.Range("N17").Value = -2669
$ws.Range("H70").Value = 1000
$ws.Range("I70").Value = 1000
$ws.Range("K70").Value = 3000
$ws.Range("M70").Value = -2685
$ws.Range("H73").Value = 1000
$ws.Range("I73").Value = 1000
$ws.Range("K73").Value = 3000
$ws.Range("M73").Value = -1908

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 18332
$ws.Range("I41").Value = 17498.5
$ws.Range("J41").Value = 19999
$ws.Range("K41").Value = 17498.5
$ws.Range("L41").Value = 19999
$ws.Range("M41").Value = -17143.5
$ws.Range("N41").Value = -20709
$ws.Range("H80").Value = 3177.8064
$ws.Range("I80").Value = 3293.625
$ws.Range("J80").Value = 3137.5217
$ws.Range("K80").Value = 3293.625
$ws.Range("L80").Value = 3137.5217
$ws.Range("M80").Value = -2295.625
$ws.Range("N80").Value = -5133.521699999999
$ws.Range("H83").Value = 3177.8064
$ws.Range("I83").Value = 3293.625
$ws.Range("J83").Value = 3137.5217
$ws.Range("K83").Value = 16468.125
$ws.Range("L83").Value = 15687.6085
$ws.Range("M83").Value = -11476.125
$ws.Range("N83").Value = -25671.6085
$ws.Range("H122").Value = 2596.6667
$ws.Range("I122").Value = 2196.7856
$ws.Range("K122").Value = 6590.3568
$ws.Range("M122").Value = -4140.3568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1517.7222
$ws.Range("I46").Value = 992.625
$ws.Range("J46").Value = 1937.8
$ws.Range("K46").Value = 992.625
$ws.Range("L46").Value = 1937.8
$ws.Range("M46").Value = -804.625
$ws.Range("N46").Value = -2313.8
$ws.Range("H82").Value = 1777
$ws.Range("I82").Value = 1774.375
$ws.Range("K82").Value = 1774.375
$ws.Range("M82").Value = -1413.375
$ws.Range("H85").Value = 1777
$ws.Range("I85").Value = 1774.375
$ws.Range("K85").Value = 1774.375
$ws.Range("M85").Value = -526.375
$ws.Range("H93").Value = 1340.1
$ws.Range("I93").Value = 1262.625
$ws.Range("J93").Value = 1650
$ws.Range("K93").Value = 1262.625
$ws.Range("L93").Value = 1650
$ws.Range("M93").Value = -14.625
$ws.Range("N93").Value = -4146
$ws.Range("H122").Value = 2977.5862
$ws.Range("I122").Value = 2410.4211
$ws.Range("K122").Value = 7231.263300000001
$ws.Range("M122").Value = -4781.263300000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4613.76
$ws.Range("I81").Value = 1577.0555
$ws.Range("K81").Value = 3154.111
$ws.Range("M81").Value = -2093.111
$ws.Range("H84").Value = 4613.76
$ws.Range("I84").Value = 1577.0555
$ws.Range("K84").Value = 15770.555
$ws.Range("M84").Value = -10466.555
$ws.Range("H100").Value = 2085.4194
$ws.Range("J100").Value = 695.8333
$ws.Range("L100").Value = 1391.6666
$ws.Range("N100").Value = -2473.6666
$ws.Range("H113").Value = 314.3889
$ws.Range("I113").Value = 210.6
$ws.Range("J113").Value = 833.3333
$ws.Range("K113").Value = 631.8
$ws.Range("L113").Value = 2499.9999
$ws.Range("M113").Value = 1538.2
$ws.Range("N113").Value = -6839.9999
